# Add release/6.0.1 to meta-sheet
#
# The sheet is a small "environment -> branch" grid:
#   row 1: dev2 | sit2 | uat2 | prod
#   row 2: release/6.0.0 | release/6.0.0 | release/6.0.0 | release/6.0.0
# A new row 3 is appended for the release/6.0.1 branch, marked with "X" in
# every environment column (not yet rolled out anywhere).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "release/6.0.1"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# New row is plain/unformatted data (not copied down from the styled
# header/data rows above), so make sure it keeps the default style.
$ws.Range("A3:D3").Style = "Normal"
